$d = $word.ActiveDocument

# 1) Fix capitalization of "bessons" -> "Bessons" in the four
#    "Dates de la campanya..." paragraphs (simple text fix, keeps
#    existing run/paragraph formatting as-is).
$d.Content.Find.Execute(
    "Constel·lació de bessons 14-23",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Constel·lació de Bessons 14-23", 2)

# 2) Rewrite the "Esteu participant..." paragraph: collapse the many
#    separate runs into a single plain run (no rPr), updating the
#    wording to reference "Constel·lació de Bessons" instead of
#    "constel·lació Perseus".
$rng = $d.Content
$found = $rng.Find.Execute("Esteu participant", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs(1)
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End

    $newText = "Esteu participant en una campanya mundial per observar i anotar la brillantor de les estrelles més febles que es poden veure, com a mitjà per mesurar la contaminació lumínica en un lloc determinat. Localitzant i observant la  Constel·lació de Bessons a la nit i comparant la brillantor de les estrelles del cel amb la brillantor que indiquen els mapes, gent de tot el món aprendran com els llums de la seva zona contribueixen a augmentar la contaminació lumínica. Les vostres aportacions a la base de dades activa faran palesa la visibilitat del cel nocturn."

    # Delete the paragraph's text (but not its trailing paragraph mark)
    # then insert the new, single run of text in its place.
    $delRng = $d.Range($pStart, $pEnd - 1)
    $delRng.Delete()
    $insRng = $d.Range($pStart, $pStart)
    $insRng.InsertAfter($newText)
}
